$d = $word.ActiveDocument

# 1. Remove the first run's text entirely:
#    "The overall summary of the abouve tabulation is highest percentages of marks "
$d.Content.Find.Execute("The overall summary of the abouve tabulation is highest percentages of marks ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 2. Replace the second run's text with the new combined sentence
$d.Content.Find.Execute("exists between Q3 to Q4.", $true, $false, $false, $false, $false, $true, 1, $false, "As Compared  from  SSC_P marks to MBA_P marks except etest_p marks all other percentages of the marks are decreased. The overall summary of the above tabulation is highest percentages of marks exists between Q3 to Q4.", 2)

# 3. Remove three now-superfluous empty paragraphs that followed that paragraph.
#    Locate the paragraph holding the new sentence, then delete the next three
#    (empty) paragraphs that immediately follow it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*As Compared*exists between Q3 to Q4.*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $afterStart = $target.Range.End
    $p1 = $target.Next()
    $p2 = $p1.Next()
    $p3 = $p2.Next()
    $p4 = $p3.Next()
    $deleteEnd = $p4.Range.Start
    $r = $d.Range($afterStart, $deleteEnd)
    $r.Delete()
}
